# "fiexed typo in list slides"
#
# On the "Index Positions Illustrated" slide, each of the seven
# "<TAB>giraffe[n] = '...'" example lines has the word "giraffe"
# swapped out for a larger, separately-sized " animal" run (80pt)
# so the line reads "<TAB> animal[n] = '...'" - matching the big
# "animal" callout used elsewhere on the slide.
#
# In every one of these paragraphs the word "giraffe" begins right
# after the leading tab character, i.e. at paragraph character
# position 2, and is always 7 characters long - so the same two
# operations (replace text, then bump font size) are applied to the
# same character range in each paragraph. Only that 7-character
# sub-range is ever touched; everything else (tabs, brackets, quote
# marks, the misspelling-flagged "i" run, the trailing tab run, etc.)
# is left completely untouched so existing run formatting survives.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraphs (1-based, within this shape) that contain a
# "<TAB>giraffe[n] = ..." line: index positions 0 through 6.
$paraIndexes = 3, 4, 5, 6, 7, 8, 9

foreach ($idx in $paraIndexes) {
  $para = $tr.Paragraphs($idx, 1)

  # "giraffe" is always the 7 characters starting right after the
  # leading tab (paragraph character #2).
  $word = $para.Characters(2, 7)
  $word.Text = " animal"

  # Re-fetch the same range (text length is unchanged: 7 chars) and
  # enlarge it to 80pt (sz="8000"), leaving the rest of the line at
  # its original size.
  $animalRun = $para.Characters(2, 7)
  $animalRun.Font.Size = 80
}
